$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.794.27"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.634.16"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "1.639.91"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "1.858.96"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "0.0₃0772"
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "25.805.00"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  -4.32%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.897"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.551"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "1.108.44"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "0.0₆0110"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
